$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at the top of the data (row 2), pushing existing data down
$ws.Rows.Item(2).Insert()

# Carry the formatting of the (now shifted) former row 2 up into the new row
$ws.Range("A3:C3").Copy($ws.Range("A2:C2"))

# Expand the table / query range to include the new row
$lo.Resize($ws.Range("A1:C24"))

# Populate the new row with the latest day's P&L data
$ws.Cells.Item(2, 1).Value = "ETHUSDT"
$ws.Cells.Item(2, 2).Value = 44587
$ws.Cells.Item(2, 3).Value = -290.8

$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$C`$24"
